$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new "Recorded By" value (System moved to front of list)
$updates = @{
    'G2' = 'System, system, backup@backdoor.com'
    'G4' = 'System, backup@backdoor.com'
    'G5' = 'System, backup@backdoor.com'
    'G8' = 'System, backup@backdoor.com'
    'G11' = 'System, dnasr281@gmail.com'
    'G17' = 'System, dnasr281@gmail.com'
    'G28' = 'System, system, backup@backdoor.com'
    'G30' = 'System, backup@backdoor.com'
    'G31' = 'System, backup@backdoor.com'
    'G34' = 'System, backup@backdoor.com'
    'G37' = 'System, dnasr281@gmail.com'
    'G43' = 'System, dnasr281@gmail.com'
    'G54' = 'System, system, backup@backdoor.com'
    'G56' = 'System, backup@backdoor.com'
    'G57' = 'System, backup@backdoor.com'
    'G60' = 'System, backup@backdoor.com'
    'G63' = 'System, dnasr281@gmail.com'
    'G69' = 'System, dnasr281@gmail.com'
    'G80' = 'System, backup@backdoor.com'
    'G81' = 'System, backup@backdoor.com'
    'G82' = 'System, backup@backdoor.com'
    'G93' = 'System, dnasr281@gmail.com'
    'G94' = 'System, dnasr281@gmail.com'
    'G96' = 'System, dnasr281@gmail.com'
    'G106' = 'System, backup@backdoor.com'
    'G107' = 'System, backup@backdoor.com'
    'G108' = 'System, backup@backdoor.com'
    'G119' = 'System, dnasr281@gmail.com'
    'G120' = 'System, dnasr281@gmail.com'
    'G122' = 'System, dnasr281@gmail.com'
    'G132' = 'System, backup@backdoor.com'
    'G133' = 'System, backup@backdoor.com'
    'G134' = 'System, backup@backdoor.com'
    'G145' = 'System, dnasr281@gmail.com'
    'G146' = 'System, dnasr281@gmail.com'
    'G148' = 'System, dnasr281@gmail.com'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

